$wb = $excel.ActiveWorkbook

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2470.074
$ws.Cells.Item(113, 10).Value = 2338.25
$ws.Cells.Item(113, 12).Value = 2338.25
$ws.Cells.Item(113, 14).Value = -8846.25

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4157.0713
$ws.Cells.Item(116, 10).Value = 3377.7778
$ws.Cells.Item(116, 12).Value = 3377.7778
$ws.Cells.Item(116, 14).Value = -10261.7778

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 2240.5186
$ws.Cells.Item(132, 9).Value = 2473.6956
$ws.Cells.Item(132, 10).Value = 899.75
$ws.Cells.Item(132, 11).Value = 7421.0868
$ws.Cells.Item(132, 12).Value = 2699.25
$ws.Cells.Item(132, 13).Value = -4891.0868
$ws.Cells.Item(132, 14).Value = -7759.25

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 2000.591
$ws.Cells.Item(141, 9).Value = 666.3143
$ws.Cells.Item(141, 11).Value = 1998.9429
$ws.Cells.Item(141, 13).Value = 3181.0571

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17732.197
$ws.Cells.Item(32, 9).Value = 20600.25
$ws.Cells.Item(32, 10).Value = 9252.739
$ws.Cells.Item(32, 11).Value = 20600.25
$ws.Cells.Item(32, 12).Value = 9252.739
$ws.Cells.Item(32, 13).Value = -20313.25
$ws.Cells.Item(32, 14).Value = -9826.739

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1143.6976
$ws.Cells.Item(61, 9).Value = 759
$ws.Cells.Item(61, 11).Value = 759
$ws.Cells.Item(61, 13).Value = -547

# ARM row 117
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1537.7435
$ws.Cells.Item(132, 9).Value = 1121.1538
$ws.Cells.Item(132, 10).Value = 2370.923
$ws.Cells.Item(132, 11).Value = 3363.4614
$ws.Cells.Item(132, 12).Value = 7112.768999999999
$ws.Cells.Item(132, 13).Value = -833.4614000000001
$ws.Cells.Item(132, 14).Value = -12172.769

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1143.6976
$ws.Cells.Item(136, 9).Value = 759
$ws.Cells.Item(136, 11).Value = 2277
$ws.Cells.Item(136, 13).Value = 273

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1519.5769
$ws.Cells.Item(99, 9).Value = 1490.5
$ws.Cells.Item(99, 10).Value = 1616.5
$ws.Cells.Item(99, 11).Value = 1490.5
$ws.Cells.Item(99, 12).Value = 1616.5
$ws.Cells.Item(99, 13).Value = 7.5
$ws.Cells.Item(99, 14).Value = -4612.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 16301.235
$ws.Cells.Item(134, 9).Value = 1262.6545
$ws.Cells.Item(134, 10).Value = 79926
$ws.Cells.Item(134, 11).Value = 3787.9635
$ws.Cells.Item(134, 12).Value = 239778
$ws.Cells.Item(134, 13).Value = -1252.9635
$ws.Cells.Item(134, 14).Value = -244848

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 148.57143
$ws.Cells.Item(19, 9).Value = 148.57143
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 148.57143
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 21.42857000000001
$ws.Cells.Item(19, 14).ClearContents()

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(24, 8).Value = 148.57143
$ws.Cells.Item(24, 9).Value = 148.57143
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 148.57143
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 21.42857000000001
$ws.Cells.Item(24, 14).ClearContents()

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2221.951
$ws.Cells.Item(31, 9).Value = 2168.372
$ws.Cells.Item(31, 11).Value = 2168.372
$ws.Cells.Item(31, 13).Value = -1873.372

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2221.951
$ws.Cells.Item(34, 9).Value = 2168.372
$ws.Cells.Item(34, 11).Value = 2168.372
$ws.Cells.Item(34, 13).Value = -1966.372

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 3221.5789
$ws.Cells.Item(99, 9).Value = 3342.4
$ws.Cells.Item(99, 11).Value = 3342.4
$ws.Cells.Item(99, 13).Value = -1844.4

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 3221.5789
$ws.Cells.Item(126, 9).Value = 3342.4
$ws.Cells.Item(126, 11).Value = 10027.2
$ws.Cells.Item(126, 13).Value = -7557.200000000001

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1693.0222
$ws.Cells.Item(132, 9).Value = 1160.4193
$ws.Cells.Item(132, 10).Value = 2872.3572
$ws.Cells.Item(132, 11).Value = 3481.2579
$ws.Cells.Item(132, 12).Value = 8617.071599999999
$ws.Cells.Item(132, 13).Value = -951.2579000000001
$ws.Cells.Item(132, 14).Value = -13677.0716

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1411.6492
$ws.Cells.Item(134, 9).Value = 1268.7174
$ws.Cells.Item(134, 10).Value = 2009.3636
$ws.Cells.Item(134, 11).Value = 3806.1522
$ws.Cells.Item(134, 12).Value = 6028.0908
$ws.Cells.Item(134, 13).Value = -1271.1522
$ws.Cells.Item(134, 14).Value = -11098.0908

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 4291.875
$ws.Cells.Item(70, 9).Value = 2947
$ws.Cells.Item(70, 10).Value = 6533.3335
$ws.Cells.Item(70, 11).Value = 8841
$ws.Cells.Item(70, 12).Value = 19600.0005
$ws.Cells.Item(70, 13).Value = -8526
$ws.Cells.Item(70, 14).Value = -20230.0005

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(73, 8).Value = 4291.875
$ws.Cells.Item(73, 9).Value = 2947
$ws.Cells.Item(73, 10).Value = 6533.3335
$ws.Cells.Item(73, 11).Value = 8841
$ws.Cells.Item(73, 12).Value = 19600.0005
$ws.Cells.Item(73, 13).Value = -7749
$ws.Cells.Item(73, 14).Value = -21784.0005

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 526959.9
$ws.Cells.Item(122, 10).Value = 715040.7
$ws.Cells.Item(122, 12).Value = 6435366.3
$ws.Cells.Item(122, 14).Value = -6440266.3

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 6123443.5
$ws.Cells.Item(131, 9).Value = 125252700
$ws.Cells.Item(131, 10).Value = 14250.68
$ws.Cells.Item(131, 11).Value = 375758100
$ws.Cells.Item(131, 12).Value = 42752.04
$ws.Cells.Item(131, 13).Value = -375753060
$ws.Cells.Item(131, 14).Value = -52832.04

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4341
$ws.Cells.Item(70, 9).Value = 4088.5
$ws.Cells.Item(70, 10).Value = 4789.8887
$ws.Cells.Item(70, 11).Value = 4088.5
$ws.Cells.Item(70, 12).Value = 4789.8887
$ws.Cells.Item(70, 13).Value = -3818.5
$ws.Cells.Item(70, 14).Value = -5329.8887

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 4341
$ws.Cells.Item(73, 9).Value = 4088.5
$ws.Cells.Item(73, 10).Value = 4789.8887
$ws.Cells.Item(73, 11).Value = 4088.5
$ws.Cells.Item(73, 12).Value = 4789.8887
$ws.Cells.Item(73, 13).Value = -3152.5
$ws.Cells.Item(73, 14).Value = -6661.8887

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2095.025
$ws.Cells.Item(132, 9).Value = 1893.7812
$ws.Cells.Item(132, 10).Value = 2900
$ws.Cells.Item(132, 11).Value = 5681.3436
$ws.Cells.Item(132, 12).Value = 8700
$ws.Cells.Item(132, 13).Value = -3151.3436
$ws.Cells.Item(132, 14).Value = -13760

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2099.3572
$ws.Cells.Item(46, 9).Value = 1741.5714
$ws.Cells.Item(46, 11).Value = 1741.5714
$ws.Cells.Item(46, 13).Value = -1553.5714

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 8788.532999999999
$ws.Cells.Item(122, 9).Value = 11982.8
$ws.Cells.Item(122, 10).Value = 2400
$ws.Cells.Item(122, 11).Value = 35948.39999999999
$ws.Cells.Item(122, 12).Value = 7200
$ws.Cells.Item(122, 13).Value = -33498.39999999999
$ws.Cells.Item(122, 14).Value = -12100

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1044.8572
$ws.Cells.Item(122, 9).Value = 1102
$ws.Cells.Item(122, 10).Value = 1022
$ws.Cells.Item(122, 11).Value = 3306
$ws.Cells.Item(122, 12).Value = 3066
$ws.Cells.Item(122, 13).Value = -856
$ws.Cells.Item(122, 14).Value = -7966

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 47500
$ws.Cells.Item(123, 10).Value = 47500
$ws.Cells.Item(123, 12).Value = 47500
$ws.Cells.Item(123, 14).Value = -57300
